$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Jenis Kelamin" column with header + values
$ws.Range("G1").Value = "Jenis Kelamin"
$ws.Range("G2").Value = "L"
$ws.Range("G3").Value = "P"
$ws.Range("G4").Value = "P"
$ws.Range("G5").Value = "L"
$ws.Range("G6").Value = "L"

# Match the header style (bold) used on row 1, columns A1:F1
$ws.Range("G1").Font.Bold = $true

# Match the left-aligned style used by the rest of the data rows (A2:F6)
$ws.Range("G2:G6").HorizontalAlignment = -4131

# Set column width to match the new narrower "Jenis Kelamin" column
$ws.Range("G1").ColumnWidth = 12.83

# Update selection / active cell as captured in the diff
$ws.Range("A1:XFD1048576").Select()
$ws.Range("H3").Activate()
